# Insert a new row at row 60. This shifts the existing rows 60-181
# down to 61-182 (all data, formatting and styles move with them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()

# Populate the newly inserted (blank) row 60 with the new record.
$ws.Cells.Item(60, 1).Value = 5
$ws.Cells.Item(60, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(60, 3).Value = "Maule"
$ws.Cells.Item(60, 4).Value = 45125
$ws.Cells.Item(60, 5).Value = 7
$ws.Cells.Item(60, 6).Value = 100112001
$ws.Cells.Item(60, 7).Value = "Berenjena"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 300
$ws.Cells.Item(60, 11).Value = 7000
$ws.Cells.Item(60, 12).Value = 7000
$ws.Cells.Item(60, 13).Value = 7000
$ws.Cells.Item(60, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(60, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 16).Value = 140
$ws.Cells.Item(60, 17).Value = 50
$ws.Cells.Item(60, 18).Value = "Hortaliza"
